$d = $word.ActiveDocument
$t = $d.Tables.Item($d.Tables.Count)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "17/07/2023"
$newRow.Cells.Item(2).Range.Text = "Thalyta Costa"
$newRow.Cells.Item(3).Range.Text = "Correção da formatação"
Write-Output "done"
